# Rename Sheet1 -> "default endpoint"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "default endpoint"

# -----------------------------------------------------------------------
# New shared-string text blocks used by rows 16 and 17 (kept as variables
# so they are typed/assigned to cells in the exact order the original
# workbook's sharedStrings table lists them: D16, B16, E16, D17, A17,
# B17, E17).
# -----------------------------------------------------------------------
$req16 = @'
{
  "score": 0,
  "is_finish": false,
  "user_id": "a3383d40-1b5c-4355-8889-1aa84b0e61f7",
  "assignment_id": "5aa3b902-4b66-4080-bcb1-63f9d9cd86cc"
}
'@

$url16 = 'https://ecos.joheee.com/googolplex/user_assignment_todo'

$resp16 = @'
{
  "status": 200,
  "message": "user_assignment_todo successfully created!",
  "data": {
    "id": "1e305472-edb8-4c0c-8006-082073ebad3a",
    "created_at": "2024-12-04T22:12:31.447Z",
    "updated_at": "2024-12-04T22:12:31.447Z",
    "score": 0,
    "is_finish": false,
    "user_id": "a3383d40-1b5c-4355-8889-1aa84b0e61f7",
    "assignment_id": "5aa3b902-4b66-4080-bcb1-63f9d9cd86cc"
  }
}
'@

$req17 = @'
{
  "user_assignment_todo_id": "1e305472-edb8-4c0c-8006-082073ebad3a",
  "answer": "this is user todo answer"
}
'@

$endpoint17 = 'user_todo_answer'

$url17 = 'https://ecos.joheee.com/googolplex/user_todo_answer'

$resp17 = @'
{
  "status": 200,
  "message": "user_todo_answer successfully created!",
  "data": {
    "id": "840785eb-dba7-47ee-9282-80ff7c12ed0a",
    "answer": "this is user todo answer",
    "created_at": "2024-12-04T22:13:32.689Z",
    "updated_at": "2024-12-04T22:13:32.689Z",
    "user_assignment_todo_id": "1e305472-edb8-4c0c-8006-082073ebad3a"
  }
}
'@

# -----------------------------------------------------------------------
# Row 16: copy formats from row 15 (same look: centered col A/C/F,
# wrap-text col D/E, hyperlink style col B), then fill in the values.
# -----------------------------------------------------------------------
$ws.Range("A15:F15").Copy()
$ws.Range("A16:F16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D16").Value = $req16
$ws.Range("B16").Value = $url16
$ws.Range("E16").Value = $resp16
$ws.Range("A16").Value = "user_assignment_todo"
$ws.Range("C16").Value = "post"
$ws.Range("F16").Value = "done"
$ws.Rows.Item(16).RowHeight = 195

[void]$ws.Hyperlinks.Add($ws.Range("B16"), $url16)
# Hyperlinks.Add() overwrites the cell style with its own built-in
# Hyperlink look; re-apply the same Hyperlink style used by B7..B15 that
# PasteSpecial already gave us, so column B keeps matching formatting.
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# -----------------------------------------------------------------------
# Row 17: same treatment.
# -----------------------------------------------------------------------
$ws.Range("A15:F15").Copy()
$ws.Range("A17:F17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D17").Value = $req17
$ws.Range("A17").Value = $endpoint17
$ws.Range("B17").Value = $url17
$ws.Range("E17").Value = $resp17
$ws.Range("C17").Value = "post"
$ws.Range("F17").Value = "done"
$ws.Rows.Item(17).RowHeight = 165

[void]$ws.Hyperlinks.Add($ws.Range("B17"), $url17)
$ws.Range("B15").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# -----------------------------------------------------------------------
# Column width tweaks.
# -----------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 32
$ws.Columns.Item(5).ColumnWidth = 83
$ws.Columns.Item(7).ColumnWidth = 29

# -----------------------------------------------------------------------
# Selection / view state: select G3, clear the old frozen top-left cell.
# -----------------------------------------------------------------------
[void]$ws.Range("G3").Select()
